$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (gender) changes from text labels (男/女) to numeric 1/0
# with a custom number format "0_ " applied.
$ws.Range("B2").Value = 1
$ws.Range("B2").NumberFormat = "0_ "

$ws.Range("B3").Value = 0
$ws.Range("B3").NumberFormat = "0_ "

# Update the sheet's current selection to E17.
$ws.Range("E17").Select()
